# Update gh-pages to output generated at 456a3b4
# Applies the same F/G column corrections to the "展览" (exhibition) sheet
# and to the "全部类型" (all-types) sheet, which mirrors it (plus one extra
# row, #23). Also bumps the single data row on "演出" (performance).

$wb = $excel.ActiveWorkbook

function Update-ExhibitionRows($ws) {
    $ws.Range("F2").Value = 320
    $ws.Range("F3").Value = 13931
    $ws.Range("F6").Value = 189
    $ws.Range("G6").Value = "不可售"
    $ws.Range("F7").Value = 289
    $ws.Range("F8").Value = 506
    $ws.Range("F10").Value = 92
    $ws.Range("F13").Value = 56
    $ws.Range("F14").Value = 466
    $ws.Range("F15").Value = 5978
    $ws.Range("F16").Value = 147
    $ws.Range("F18").Value = 991
    $ws.Range("F19").Value = 152
    $ws.Range("F20").Value = 65
    $ws.Range("F21").Value = 163
    $ws.Range("F22").Value = 317
}

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
Update-ExhibitionRows $wsExhibit

# Sheet "演出"
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 16

# Sheet "全部类型" (same rows as 展览, plus row 23)
$wsAll = $wb.Worksheets.Item("全部类型")
Update-ExhibitionRows $wsAll
$wsAll.Range("F23").Value = 16
